# Fix random assignation of metadata to files
# Updates two description cells on the "Attribute description" sheet:
#  - isolation_source description gains an "unknown" option
#  - collection_date description gains an "unknown" option, with the
#    allowed format tokens (and "unknown") shown in bold

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Attribute description")

# --- isolation_source (row 19) description: mention "unknown" too ---
$ws.Range("C19").Value = "Information about the isolation source (i.e. blood, laboratory experiment, urine, unknown...)"

# --- collection_date (row 24) description: add "or unknown" option, bold the keywords ---
$dateCell = $ws.Range("C24")
$dateCell.Value = "The date of the sample collection. Use one of the following format: YYYY-MM-DD, YYYY-MM or YYYY or unknown."

# Bold the four format keywords
$dateCell.Characters(69, 10).Font.Bold = $true   # YYYY-MM-DD
$dateCell.Characters(81, 7).Font.Bold = $true    # YYYY-MM
$dateCell.Characters(92, 4).Font.Bold = $true    # YYYY
$dateCell.Characters(100, 7).Font.Bold = $true   # unknown

# Touch the font size on the connecting (non-bold) runs so every run but the
# leading sentence carries explicit run-level formatting, matching how Excel
# records mixed bold/non-bold runs inside one cell.
$dateCell.Characters(79, 2).Font.Size = 16    # ", "
$dateCell.Characters(88, 4).Font.Size = 16    # " or "
$dateCell.Characters(96, 4).Font.Size = 16    # " or "
$dateCell.Characters(107, 1).Font.Size = 16   # "."

# A new (empty) column D shows up next to the edited description, sized to
# match the author's workbook, along with a bold font carried over from the
# edit above.
$ws.Columns("D").ColumnWidth = 28.375
$ws.Range("D24").Font.Bold = $true

# Restore the sheet selection/scroll state left after the edit
$ws.Activate()
$ws.Range("C26").Select()
